# Apply commit "change StockTable to StockHistory":
#  - rename the "在庫" (stock table) sheet to "在庫履歴" (stock history)
#  - rework its field-definition rows: the sheet now mirrors the 原価
#    (cost master) sheet's 在庫* columns instead of 工場 (factory) ID alone,
#    and gains 品目ID (item id), 時刻 (time) and 在庫数量 (stock quantity)
#    rows, with the stock-amount row pushed down and re-labelled.
#  - tweak the active-cell selections left behind on the 原価 and
#    在庫履歴 sheets.

$wb = $excel.ActiveWorkbook

# --- rename the sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("在庫")
$ws.Name = "在庫履歴"

# --- title row ----------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "在庫履歴"
$ws.Cells.Item(1, 2).Value = "stock_history"

# --- make room: insert 4 new field rows before the old "stock amount"
#     row (old row 5), shifting everything from old row 5 onward down to
#     row 9, and the trailing styled blank rows from old 7..22 to 11..26.
$ws.Rows("5:8").Insert()

# --- row 4: 工場ID / s_f_id now sourced from 原価 instead of 工場 --------
$ws.Cells.Item(4, 1).Value = "工場ID"
$ws.Cells.Item(4, 2).Value = "s_f_id"
$ws.Cells.Item(4, 3).Formula = "=原価!C3"
$ws.Cells.Item(4, 4).Formula = "=原価!D3"

# --- row 5 (new): 品目ID / s_i_id ---------------------------------------
$ws.Cells.Item(5, 1).Value = "品目ID"
$ws.Cells.Item(5, 2).Value = "s_i_id"
$ws.Cells.Item(5, 2).WrapText = $true
$ws.Cells.Item(5, 3).Formula = "=原価!C4"
$ws.Cells.Item(5, 4).Formula = "=原価!D4"
$ws.Cells.Item(5, 6).Value = "PK"

# --- row 6 (new): 時刻 / s_time ------------------------------------------
$ws.Cells.Item(6, 1).Value = "時刻"
$ws.Cells.Item(6, 2).Value = "s_time"
$ws.Cells.Item(6, 2).WrapText = $true
$ws.Cells.Item(6, 3).Value = "time"
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 6).Value = "PK"
$ws.Cells.Item(6, 7).Value = "時分秒"

# --- row 7: 在庫単位 / s_stock_unit (not a key field -> no wrap style) --
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 1).Value = "在庫単位"
$ws.Cells.Item(7, 2).Value = "s_stock_unit"
$ws.Cells.Item(7, 3).Formula = "=原価!C5"
$ws.Cells.Item(7, 4).Formula = "=原価!D5"
$ws.Cells.Item(7, 7).Value = "原価マスター．在庫単位"

# --- row 8 (new): 在庫数量 / s_stock_quantity ----------------------------
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 1).Value = "在庫数量"
$ws.Cells.Item(8, 2).Value = "s_stock_quantity"
$ws.Cells.Item(8, 3).Formula = "=原価!C6"
$ws.Cells.Item(8, 4).Formula = "=原価!D6"
$ws.Cells.Item(8, 5).Formula = "=原価!E6"
$ws.Cells.Item(8, 7).Value = "原価マスター．在庫数量"

# --- row 9 (was row 5): 在庫金額 / s_stock_amount, re-labelled ----------
$ws.Cells.Item(9, 2).ClearFormats()
$ws.Cells.Item(9, 1).Value = "在庫金額"
$ws.Cells.Item(9, 2).Value = "s_stock_amount"
$ws.Cells.Item(9, 3).Formula = "=原価!C7"
$ws.Cells.Item(9, 4).Formula = "=原価!D7"
$ws.Cells.Item(9, 5).Formula = "=原価!E7"
$ws.Cells.Item(9, 7).Value = "原価マスター．在庫金額"

# --- selections left behind by the author ------------------------------
$ws.Range("E12").Select()

$costSheet = $wb.Worksheets.Item("原価")
$costSheet.Range("D10").Select()

$ws.Activate()
